# Auto-generated edit script: rebuild the NIP list in column A
# (rows 2-1161) to match the corrected master list, per commit
# "develop-payroll:fix nip in template for penghasilan tidak teratur".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastNewRow = 1161

$data = New-Object "object[,]" 1160,1
$data[0,0] = "'00110"
$data[1,0] = "'00112"
$data[2,0] = "'00113"
$data[3,0] = "'00114"
$data[4,0] = "'00115"
$data[5,0] = "'00116"
$data[6,0] = "'00117"
$data[7,0] = "'00118"
$data[8,0] = "'00119"
$data[9,0] = "'00120"
$data[10,0] = "'00121"
$data[11,0] = "'00122"
$data[12,0] = "'00126"
$data[13,0] = "'00137"
$data[14,0] = "'00161"
$data[15,0] = "'00163"
$data[16,0] = "'00164"
$data[17,0] = "'00166"
$data[18,0] = "'00168"
$data[19,0] = "'00171"
$data[20,0] = "'00173"
$data[21,0] = "'00175"
$data[22,0] = "'00191"
$data[23,0] = "'00193"
$data[24,0] = "'00196"
$data[25,0] = "'00198"
$data[26,0] = "'00201"
$data[27,0] = "'00202"
$data[28,0] = "'00203"
$data[29,0] = "'00204"
$data[30,0] = "'00208"
$data[31,0] = "'00210"
$data[32,0] = "'00218"
$data[33,0] = "'00227"
$data[34,0] = "'00230"
$data[35,0] = "'00231"
$data[36,0] = "'00234"
$data[37,0] = "'00236"
$data[38,0] = "'00245"
$data[39,0] = "'00246"
$data[40,0] = "'00247"
$data[41,0] = "'00248"
$data[42,0] = "'00251"
$data[43,0] = "'00252"
$data[44,0] = "'00253"
$data[45,0] = "'00254"
$data[46,0] = "'00256"
$data[47,0] = "'00257"
$data[48,0] = "'00258"
$data[49,0] = "'00259"
$data[50,0] = "'00265"
$data[51,0] = "'00269"
$data[52,0] = "'00270"
$data[53,0] = "'00272"
$data[54,0] = "'00273"
$data[55,0] = "'00274"
$data[56,0] = "'00277"
$data[57,0] = "'00279"
$data[58,0] = "'00280"
$data[59,0] = "'00282"
$data[60,0] = "'00284"
$data[61,0] = "'00286"
$data[62,0] = "'00289"
$data[63,0] = "'00294"
$data[64,0] = "'00295"
$data[65,0] = "'00296"
$data[66,0] = "'00299"
$data[67,0] = "'00304"
$data[68,0] = "'00312"
$data[69,0] = "'00313"
$data[70,0] = "'00314"
$data[71,0] = "'00316"
$data[72,0] = "'00318"
$data[73,0] = "'00320"
$data[74,0] = "'00324"
$data[75,0] = "'00325"
$data[76,0] = "'00326"
$data[77,0] = "'00330"
$data[78,0] = "'00331"
$data[79,0] = "'00334"
$data[80,0] = "'00335"
$data[81,0] = "'00336"
$data[82,0] = "'00337"
$data[83,0] = "'00343"
$data[84,0] = "'00344"
$data[85,0] = "'00345"
$data[86,0] = "'00346"
$data[87,0] = "'00349"
$data[88,0] = "'00353"
$data[89,0] = "'00355"
$data[90,0] = "'00356"
$data[91,0] = "'00358"
$data[92,0] = "'00360"
$data[93,0] = "'00361"
$data[94,0] = "'00362"
$data[95,0] = "'00363"
$data[96,0] = "'00364"
$data[97,0] = "'00365"
$data[98,0] = "'00376"
$data[99,0] = "'00379"
$data[100,0] = "'00380"
$data[101,0] = "'00381"
$data[102,0] = "'00385"
$data[103,0] = "'00386"
$data[104,0] = "'00387"
$data[105,0] = "'00390"
$data[106,0] = "'00394"
$data[107,0] = "'00408"
$data[108,0] = "'00409"
$data[109,0] = "'00416"
$data[110,0] = "'00417"
$data[111,0] = "'00420"
$data[112,0] = "'00421"
$data[113,0] = "'00422"
$data[114,0] = "'00423"
$data[115,0] = "'00424"
$data[116,0] = "'00425"
$data[117,0] = "'00436"
$data[118,0] = "'00437"
$data[119,0] = "'00438"
$data[120,0] = "'00439"
$data[121,0] = "'00440"
$data[122,0] = "'00441"
$data[123,0] = "'00442"
$data[124,0] = "'00443"
$data[125,0] = "'00447"
$data[126,0] = "'00448"
$data[127,0] = "'00449"
$data[128,0] = "'00451"
$data[129,0] = "'00452"
$data[130,0] = "'00455"
$data[131,0] = "'00456"
$data[132,0] = "'00457"
$data[133,0] = "'00458"
$data[134,0] = "'00459"
$data[135,0] = "'00464"
$data[136,0] = "'00465"
$data[137,0] = "'00466"
$data[138,0] = "'00467"
$data[139,0] = "'00468"
$data[140,0] = "'00469"
$data[141,0] = "'00470"
$data[142,0] = "'00471"
$data[143,0] = "'00472"
$data[144,0] = "'00474"
$data[145,0] = "'00475"
$data[146,0] = "'00477"
$data[147,0] = "'00479"
$data[148,0] = "'00480"
$data[149,0] = "'00488"
$data[150,0] = "'00491"
$data[151,0] = "'00492"
$data[152,0] = "'00495"
$data[153,0] = "'00496"
$data[154,0] = "'00497"
$data[155,0] = "'00499"
$data[156,0] = "'00500"
$data[157,0] = "'00503"
$data[158,0] = "'00504"
$data[159,0] = "'00505"
$data[160,0] = "'00507"
$data[161,0] = "'00508"
$data[162,0] = "'00509"
$data[163,0] = "'00511"
$data[164,0] = "'00512"
$data[165,0] = "'00513"
$data[166,0] = "'00514"
$data[167,0] = "'00515"
$data[168,0] = "'00517"
$data[169,0] = "'00519"
$data[170,0] = "'00520"
$data[171,0] = "'00521"
$data[172,0] = "'00522"
$data[173,0] = "'00524"
$data[174,0] = "'00525"
$data[175,0] = "'00527"
$data[176,0] = "'00530"
$data[177,0] = "'00532"
$data[178,0] = "'00534"
$data[179,0] = "'00535"
$data[180,0] = "'00536"
$data[181,0] = "'00537"
$data[182,0] = "'00538"
$data[183,0] = "'00541"
$data[184,0] = "'00542"
$data[185,0] = "'00545"
$data[186,0] = "'00547"
$data[187,0] = "'00548"
$data[188,0] = "'00549"
$data[189,0] = "'00550"
$data[190,0] = "'00551"
$data[191,0] = "'00558"
$data[192,0] = "'00559"
$data[193,0] = "'00560"
$data[194,0] = "'00561"
$data[195,0] = "'00562"
$data[196,0] = "'00563"
$data[197,0] = "'00564"
$data[198,0] = "'00565"
$data[199,0] = "'00566"
$data[200,0] = "'00567"
$data[201,0] = "'00568"
$data[202,0] = "'00569"
$data[203,0] = "'00570"
$data[204,0] = "'00571"
$data[205,0] = "'00572"
$data[206,0] = "'00573"
$data[207,0] = "'00574"
$data[208,0] = "'00575"
$data[209,0] = "'00577"
$data[210,0] = "'00578"
$data[211,0] = "'00579"
$data[212,0] = "'00581"
$data[213,0] = "'00582"
$data[214,0] = "'00584"
$data[215,0] = "'00585"
$data[216,0] = "'00590"
$data[217,0] = "'00592"
$data[218,0] = "'00594"
$data[219,0] = "'00595"
$data[220,0] = "'00597"
$data[221,0] = "'00598"
$data[222,0] = "'00599"
$data[223,0] = "'00600"
$data[224,0] = "'00601"
$data[225,0] = "'00602"
$data[226,0] = "'00603"
$data[227,0] = "'00604"
$data[228,0] = "'00605"
$data[229,0] = "'00607"
$data[230,0] = "'00608"
$data[231,0] = "'00609"
$data[232,0] = "'00610"
$data[233,0] = "'00611"
$data[234,0] = "'00612"
$data[235,0] = "'00613"
$data[236,0] = "'00614"
$data[237,0] = "'00615"
$data[238,0] = "'00618"
$data[239,0] = "'00619"
$data[240,0] = "'00620"
$data[241,0] = "'00621"
$data[242,0] = "'00624"
$data[243,0] = "'00625"
$data[244,0] = "'00626"
$data[245,0] = "'00627"
$data[246,0] = "'00628"
$data[247,0] = "'00629"
$data[248,0] = "'00630"
$data[249,0] = "'00632"
$data[250,0] = "'00633"
$data[251,0] = "'00635"
$data[252,0] = "'00636"
$data[253,0] = "'00637"
$data[254,0] = "'00638"
$data[255,0] = "'00639"
$data[256,0] = "'00640"
$data[257,0] = "'00641"
$data[258,0] = "'00643"
$data[259,0] = "'00645"
$data[260,0] = "'00646"
$data[261,0] = "'00650"
$data[262,0] = "'00651"
$data[263,0] = "'00652"
$data[264,0] = "'00656"
$data[265,0] = "'00657"
$data[266,0] = "'00658"
$data[267,0] = "'00659"
$data[268,0] = "'00660"
$data[269,0] = "'00661"
$data[270,0] = "'00662"
$data[271,0] = "'00663"
$data[272,0] = "'00664"
$data[273,0] = "'00666"
$data[274,0] = "'00667"
$data[275,0] = "'00668"
$data[276,0] = "'00670"
$data[277,0] = "'00671"
$data[278,0] = "'00672"
$data[279,0] = "'00674"
$data[280,0] = "'00675"
$data[281,0] = "'00676"
$data[282,0] = "'00677"
$data[283,0] = "'00678"
$data[284,0] = "'00680"
$data[285,0] = "'00681"
$data[286,0] = "'00684"
$data[287,0] = "'00686"
$data[288,0] = "'00687"
$data[289,0] = "'00688"
$data[290,0] = "'00689"
$data[291,0] = "'00693"
$data[292,0] = "'00694"
$data[293,0] = "'00698"
$data[294,0] = "'00700"
$data[295,0] = "'00701"
$data[296,0] = "'00702"
$data[297,0] = "'00703"
$data[298,0] = "'00704"
$data[299,0] = "'00705"
$data[300,0] = "'00707"
$data[301,0] = "'00709"
$data[302,0] = "'00711"
$data[303,0] = "'00713"
$data[304,0] = "'00715"
$data[305,0] = "'00716"
$data[306,0] = "'00721"
$data[307,0] = "'00724"
$data[308,0] = "'00726"
$data[309,0] = "'00727"
$data[310,0] = "'00728"
$data[311,0] = "'00729"
$data[312,0] = "'00730"
$data[313,0] = "'00731"
$data[314,0] = "'00732"
$data[315,0] = "'00733"
$data[316,0] = "'00735"
$data[317,0] = "'00736"
$data[318,0] = "'00738"
$data[319,0] = "'00739"
$data[320,0] = "'00740"
$data[321,0] = "'00743"
$data[322,0] = "'00744"
$data[323,0] = "'00745"
$data[324,0] = "'00746"
$data[325,0] = "'00747"
$data[326,0] = "'00748"
$data[327,0] = "'00749"
$data[328,0] = "'00753"
$data[329,0] = "'00754"
$data[330,0] = "'00755"
$data[331,0] = "'00757"
$data[332,0] = "'00758"
$data[333,0] = "'00759"
$data[334,0] = "'00761"
$data[335,0] = "'00763"
$data[336,0] = "'00764"
$data[337,0] = "'00765"
$data[338,0] = "'00766"
$data[339,0] = "'00767"
$data[340,0] = "'00769"
$data[341,0] = "'00771"
$data[342,0] = "'00772"
$data[343,0] = "'00773"
$data[344,0] = "'00775"
$data[345,0] = "'00777"
$data[346,0] = "'00779"
$data[347,0] = "'00782"
$data[348,0] = "'00783"
$data[349,0] = "'00785"
$data[350,0] = "'00786"
$data[351,0] = "'00787"
$data[352,0] = "'00789"
$data[353,0] = "'00790"
$data[354,0] = "'00793"
$data[355,0] = "'00795"
$data[356,0] = "'00796"
$data[357,0] = "'00797"
$data[358,0] = "'00800"
$data[359,0] = "'00801"
$data[360,0] = "'00802"
$data[361,0] = "'00803"
$data[362,0] = "'00805"
$data[363,0] = "'00806"
$data[364,0] = "'00807"
$data[365,0] = "'00808"
$data[366,0] = "'00809"
$data[367,0] = "'00810"
$data[368,0] = "'00811"
$data[369,0] = "'00813"
$data[370,0] = "'00818"
$data[371,0] = "'00819"
$data[372,0] = "'00820"
$data[373,0] = "'00821"
$data[374,0] = "'00823"
$data[375,0] = "'00825"
$data[376,0] = "'00827"
$data[377,0] = "'00829"
$data[378,0] = "'00830"
$data[379,0] = "'00831"
$data[380,0] = "'00833"
$data[381,0] = "'00835"
$data[382,0] = "'00836"
$data[383,0] = "'00837"
$data[384,0] = "'00838"
$data[385,0] = "'00839"
$data[386,0] = "'00842"
$data[387,0] = "'00843"
$data[388,0] = "'00845"
$data[389,0] = "'00847"
$data[390,0] = "'00848"
$data[391,0] = "'00849"
$data[392,0] = "'00850"
$data[393,0] = "'00852"
$data[394,0] = "'00855"
$data[395,0] = "'00856"
$data[396,0] = "'00858"
$data[397,0] = "'00859"
$data[398,0] = "'00862"
$data[399,0] = "'00863"
$data[400,0] = "'00865"
$data[401,0] = "'00866"
$data[402,0] = "'00867"
$data[403,0] = "'00874"
$data[404,0] = "'00875"
$data[405,0] = "'00876"
$data[406,0] = "'00877"
$data[407,0] = "'00878"
$data[408,0] = "'00879"
$data[409,0] = "'00880"
$data[410,0] = "'00882"
$data[411,0] = "'00883"
$data[412,0] = "'00891"
$data[413,0] = "'00894"
$data[414,0] = "'00898"
$data[415,0] = "'00899"
$data[416,0] = "'00901"
$data[417,0] = "'00903"
$data[418,0] = "'00905"
$data[419,0] = "'00906"
$data[420,0] = "'00910"
$data[421,0] = "'00912"
$data[422,0] = "'00914"
$data[423,0] = "'00917"
$data[424,0] = "'00919"
$data[425,0] = "'00921"
$data[426,0] = "'00923"
$data[427,0] = "'00926"
$data[428,0] = "'00929"
$data[429,0] = "'00931"
$data[430,0] = "'00932"
$data[431,0] = "'00933"
$data[432,0] = "'00934"
$data[433,0] = "'00936"
$data[434,0] = "'00938"
$data[435,0] = "'00939"
$data[436,0] = "'00941"
$data[437,0] = "'00942"
$data[438,0] = "'00946"
$data[439,0] = "'00947"
$data[440,0] = "'00949"
$data[441,0] = "'00950"
$data[442,0] = "'00951"
$data[443,0] = "'00952"
$data[444,0] = "'00953"
$data[445,0] = "'00954"
$data[446,0] = "'00956"
$data[447,0] = "'00957"
$data[448,0] = "'00958"
$data[449,0] = "'00960"
$data[450,0] = "'00961"
$data[451,0] = "'00962"
$data[452,0] = "'00965"
$data[453,0] = "'00966"
$data[454,0] = "'00967"
$data[455,0] = "'00968"
$data[456,0] = "'00969"
$data[457,0] = "'00971"
$data[458,0] = "'00972"
$data[459,0] = "'00973"
$data[460,0] = "'00975"
$data[461,0] = "'00976"
$data[462,0] = "'00977"
$data[463,0] = "'00978"
$data[464,0] = "'00979"
$data[465,0] = "'00980"
$data[466,0] = "'00981"
$data[467,0] = "'00982"
$data[468,0] = "'00983"
$data[469,0] = "'00984"
$data[470,0] = "'00985"
$data[471,0] = "'00986"
$data[472,0] = "'00987"
$data[473,0] = "'00988"
$data[474,0] = "'00990"
$data[475,0] = "'00991"
$data[476,0] = "'00994"
$data[477,0] = "'00995"
$data[478,0] = "'00996"
$data[479,0] = "'00999"
$data[480,0] = "'01000"
$data[481,0] = "'01002"
$data[482,0] = "'01003"
$data[483,0] = "'01004"
$data[484,0] = "'01007"
$data[485,0] = "'01009"
$data[486,0] = "'01011"
$data[487,0] = "'01012"
$data[488,0] = "'01014"
$data[489,0] = "'01016"
$data[490,0] = "'01018"
$data[491,0] = "'01019"
$data[492,0] = "'01020"
$data[493,0] = "'01021"
$data[494,0] = "'01023"
$data[495,0] = "'01024"
$data[496,0] = "'01025"
$data[497,0] = "'01027"
$data[498,0] = "'01029"
$data[499,0] = "'01030"
$data[500,0] = "'01031"
$data[501,0] = "'01033"
$data[502,0] = "'01034"
$data[503,0] = "'01035"
$data[504,0] = "'01036"
$data[505,0] = "'01037"
$data[506,0] = "'01038"
$data[507,0] = "'01041"
$data[508,0] = "'01043"
$data[509,0] = "'01044"
$data[510,0] = "'01045"
$data[511,0] = "'01046"
$data[512,0] = "'01047"
$data[513,0] = "'01048"
$data[514,0] = "'01049"
$data[515,0] = "'01050"
$data[516,0] = "'01051"
$data[517,0] = "'01053"
$data[518,0] = "'01054"
$data[519,0] = "'01055"
$data[520,0] = "'01056"
$data[521,0] = "'01057"
$data[522,0] = "'01058"
$data[523,0] = "'01059"
$data[524,0] = "'01060"
$data[525,0] = "'01061"
$data[526,0] = "'01062"
$data[527,0] = "'01063"
$data[528,0] = "'01064"
$data[529,0] = "'01066"
$data[530,0] = "'01067"
$data[531,0] = "'01068"
$data[532,0] = "'01069"
$data[533,0] = "'01070"
$data[534,0] = "'01071"
$data[535,0] = "'01072"
$data[536,0] = "'01073"
$data[537,0] = "'01074"
$data[538,0] = "'01075"
$data[539,0] = "'01076"
$data[540,0] = "'01077"
$data[541,0] = "'01079"
$data[542,0] = "'01080"
$data[543,0] = "'01082"
$data[544,0] = "'01084"
$data[545,0] = "'01085"
$data[546,0] = "'01086"
$data[547,0] = "'01087"
$data[548,0] = "'01088"
$data[549,0] = "'01091"
$data[550,0] = "'01092"
$data[551,0] = "'01094"
$data[552,0] = "'01095"
$data[553,0] = "'01096"
$data[554,0] = "'01097"
$data[555,0] = "'01099"
$data[556,0] = "'01100"
$data[557,0] = "'01102"
$data[558,0] = "'01103"
$data[559,0] = "'01105"
$data[560,0] = "'01106"
$data[561,0] = "'01107"
$data[562,0] = "'01108"
$data[563,0] = "'01110"
$data[564,0] = "'01111"
$data[565,0] = "'01113"
$data[566,0] = "'01114"
$data[567,0] = "'01115"
$data[568,0] = "'01116"
$data[569,0] = "'01117"
$data[570,0] = "'01118"
$data[571,0] = "'01119"
$data[572,0] = "'01120"
$data[573,0] = "'01121"
$data[574,0] = "'01123"
$data[575,0] = "'01124"
$data[576,0] = "'01126"
$data[577,0] = "'01127"
$data[578,0] = "'01130"
$data[579,0] = "'01132"
$data[580,0] = "'01133"
$data[581,0] = "'01135"
$data[582,0] = "'01138"
$data[583,0] = "'01140"
$data[584,0] = "'01142"
$data[585,0] = "'01143"
$data[586,0] = "'01144"
$data[587,0] = "'01145"
$data[588,0] = "'01146"
$data[589,0] = "'01148"
$data[590,0] = "'01150"
$data[591,0] = "'01151"
$data[592,0] = "'01152"
$data[593,0] = "'01153"
$data[594,0] = "'01154"
$data[595,0] = "'01155"
$data[596,0] = "'01157"
$data[597,0] = "'01158"
$data[598,0] = "'01159"
$data[599,0] = "'01160"
$data[600,0] = "'01161"
$data[601,0] = "'01162"
$data[602,0] = "'01163"
$data[603,0] = "'01164"
$data[604,0] = "'01165"
$data[605,0] = "'01166"
$data[606,0] = "'01167"
$data[607,0] = "'01168"
$data[608,0] = "'01169"
$data[609,0] = "'01170"
$data[610,0] = "'01171"
$data[611,0] = "'01172"
$data[612,0] = "'01173"
$data[613,0] = "'01174"
$data[614,0] = "'01175"
$data[615,0] = "'01176"
$data[616,0] = "'01177"
$data[617,0] = "'01179"
$data[618,0] = "'01181"
$data[619,0] = "'01183"
$data[620,0] = "'01185"
$data[621,0] = "'01186"
$data[622,0] = "'01187"
$data[623,0] = "'01188"
$data[624,0] = "'01189"
$data[625,0] = "'01190"
$data[626,0] = "'01192"
$data[627,0] = "'01194"
$data[628,0] = "'01195"
$data[629,0] = "'01196"
$data[630,0] = "'01197"
$data[631,0] = "'01198"
$data[632,0] = "'01199"
$data[633,0] = "'01200"
$data[634,0] = "'01201"
$data[635,0] = "'01202"
$data[636,0] = "'01204"
$data[637,0] = "'01205"
$data[638,0] = "'01207"
$data[639,0] = "'01208"
$data[640,0] = "'01209"
$data[641,0] = "'01210"
$data[642,0] = "'01211"
$data[643,0] = "'01212"
$data[644,0] = "'01213"
$data[645,0] = "'01214"
$data[646,0] = "'01216"
$data[647,0] = "'01217"
$data[648,0] = "'01218"
$data[649,0] = "'01219"
$data[650,0] = "'01220"
$data[651,0] = "'01221"
$data[652,0] = "'01222"
$data[653,0] = "'01223"
$data[654,0] = "'01224"
$data[655,0] = "'01225"
$data[656,0] = "'01226"
$data[657,0] = "'01227"
$data[658,0] = "'01231"
$data[659,0] = "'01232"
$data[660,0] = "'01233"
$data[661,0] = "'01235"
$data[662,0] = "'01236"
$data[663,0] = "'01237"
$data[664,0] = "'01238"
$data[665,0] = "'01239"
$data[666,0] = "'01240"
$data[667,0] = "'01241"
$data[668,0] = "'01243"
$data[669,0] = "'01244"
$data[670,0] = "'01245"
$data[671,0] = "'01246"
$data[672,0] = "'01247"
$data[673,0] = "'01248"
$data[674,0] = "'01250"
$data[675,0] = "'01251"
$data[676,0] = "'01252"
$data[677,0] = "'01253"
$data[678,0] = "'01255"
$data[679,0] = "'01256"
$data[680,0] = "'01257"
$data[681,0] = "'01259"
$data[682,0] = "'01260"
$data[683,0] = "'01261"
$data[684,0] = "'01263"
$data[685,0] = "'01264"
$data[686,0] = "'01265"
$data[687,0] = "'01266"
$data[688,0] = "'01267"
$data[689,0] = "'01269"
$data[690,0] = "'01271"
$data[691,0] = "'01273"
$data[692,0] = "'01275"
$data[693,0] = "'01276"
$data[694,0] = "'01277"
$data[695,0] = "'01278"
$data[696,0] = "'01279"
$data[697,0] = "'01280"
$data[698,0] = "'01281"
$data[699,0] = "'01282"
$data[700,0] = "'01283"
$data[701,0] = "'01284"
$data[702,0] = "'01285"
$data[703,0] = "'01286"
$data[704,0] = "'01288"
$data[705,0] = "'01289"
$data[706,0] = "'01290"
$data[707,0] = "'01291"
$data[708,0] = "'01292"
$data[709,0] = "'01293"
$data[710,0] = "'01294"
$data[711,0] = "'01295"
$data[712,0] = "'01296"
$data[713,0] = "'01298"
$data[714,0] = "'01300"
$data[715,0] = "'01301"
$data[716,0] = "'01302"
$data[717,0] = "'01303"
$data[718,0] = "'01304"
$data[719,0] = "'01305"
$data[720,0] = "'01306"
$data[721,0] = "'01307"
$data[722,0] = "'01309"
$data[723,0] = "'01310"
$data[724,0] = "'01312"
$data[725,0] = "'01313"
$data[726,0] = "'01314"
$data[727,0] = "'01315"
$data[728,0] = "'01316"
$data[729,0] = "'01317"
$data[730,0] = "'01318"
$data[731,0] = "'01319"
$data[732,0] = "'01320"
$data[733,0] = "'01321"
$data[734,0] = "'01322"
$data[735,0] = "'01323"
$data[736,0] = "'01324"
$data[737,0] = "'01325"
$data[738,0] = "'01326"
$data[739,0] = "'01327"
$data[740,0] = "'01328"
$data[741,0] = "'01329"
$data[742,0] = "'01330"
$data[743,0] = "'01331"
$data[744,0] = "'01332"
$data[745,0] = "'01333"
$data[746,0] = "'01335"
$data[747,0] = "'01336"
$data[748,0] = "'01337"
$data[749,0] = "'01338"
$data[750,0] = "'01339"
$data[751,0] = "'01340"
$data[752,0] = "'01341"
$data[753,0] = "'01342"
$data[754,0] = "'01344"
$data[755,0] = "'01345"
$data[756,0] = "'01346"
$data[757,0] = "'01347"
$data[758,0] = "'01348"
$data[759,0] = "'01349"
$data[760,0] = "'01350"
$data[761,0] = "'01351"
$data[762,0] = "'01352"
$data[763,0] = "'01353"
$data[764,0] = "'01355"
$data[765,0] = "'01356"
$data[766,0] = "'01357"
$data[767,0] = "'01358"
$data[768,0] = "'01359"
$data[769,0] = "'01360"
$data[770,0] = "'01361"
$data[771,0] = "'01362"
$data[772,0] = "'01363"
$data[773,0] = "'01364"
$data[774,0] = "'01365"
$data[775,0] = "'01366"
$data[776,0] = "'01367"
$data[777,0] = "'01368"
$data[778,0] = "'01369"
$data[779,0] = "'01370"
$data[780,0] = "'01371"
$data[781,0] = "'01373"
$data[782,0] = "'01375"
$data[783,0] = "'01376"
$data[784,0] = "'01377"
$data[785,0] = "'01378"
$data[786,0] = "'01379"
$data[787,0] = "'01380"
$data[788,0] = "'01383"
$data[789,0] = "'01384"
$data[790,0] = "'01385"
$data[791,0] = "'01386"
$data[792,0] = "'01387"
$data[793,0] = "'01388"
$data[794,0] = "'01389"
$data[795,0] = "'01391"
$data[796,0] = "'01392"
$data[797,0] = "'01393"
$data[798,0] = "'01394"
$data[799,0] = "'01395"
$data[800,0] = "'01396"
$data[801,0] = "'01397"
$data[802,0] = "'01398"
$data[803,0] = "'01399"
$data[804,0] = "'01400"
$data[805,0] = "'01401"
$data[806,0] = "'01402"
$data[807,0] = "'01403"
$data[808,0] = "'01405"
$data[809,0] = "'01406"
$data[810,0] = "'01407"
$data[811,0] = "'01408"
$data[812,0] = "'01410"
$data[813,0] = "'01411"
$data[814,0] = "'01412"
$data[815,0] = "'01413"
$data[816,0] = "'01414"
$data[817,0] = "'01415"
$data[818,0] = "'01416"
$data[819,0] = "'01417"
$data[820,0] = "'01418"
$data[821,0] = "'01419"
$data[822,0] = "'01420"
$data[823,0] = "'01421"
$data[824,0] = "'01422"
$data[825,0] = "'01423"
$data[826,0] = "'01424"
$data[827,0] = "'01425"
$data[828,0] = "'01426"
$data[829,0] = "'01427"
$data[830,0] = "'01428"
$data[831,0] = "'01429"
$data[832,0] = "'01430"
$data[833,0] = "'01431"
$data[834,0] = "'01432"
$data[835,0] = "'01433"
$data[836,0] = "'01434"
$data[837,0] = "'01436"
$data[838,0] = "'01437"
$data[839,0] = "'01438"
$data[840,0] = "'01439"
$data[841,0] = "'01440"
$data[842,0] = "'01441"
$data[843,0] = "'01442"
$data[844,0] = "'01443"
$data[845,0] = "'01444"
$data[846,0] = "'01445"
$data[847,0] = "'01446"
$data[848,0] = "'01447"
$data[849,0] = "'01448"
$data[850,0] = "'01449"
$data[851,0] = "'01450"
$data[852,0] = "'01451"
$data[853,0] = "'01452"
$data[854,0] = "'01453"
$data[855,0] = "'01454"
$data[856,0] = "'01455"
$data[857,0] = "'01456"
$data[858,0] = "'01457"
$data[859,0] = "'01458"
$data[860,0] = "'01459"
$data[861,0] = "'01460"
$data[862,0] = "'01461"
$data[863,0] = "'01462"
$data[864,0] = "'01463"
$data[865,0] = "'01464"
$data[866,0] = "'01465"
$data[867,0] = "'01466"
$data[868,0] = "'01467"
$data[869,0] = "'01468"
$data[870,0] = "'01469"
$data[871,0] = "'01470"
$data[872,0] = "'01471"
$data[873,0] = "'01472"
$data[874,0] = "'01473"
$data[875,0] = "'01474"
$data[876,0] = "'01475"
$data[877,0] = "'01476"
$data[878,0] = "'01478"
$data[879,0] = "'01479"
$data[880,0] = "'01480"
$data[881,0] = "'01481"
$data[882,0] = "'01482"
$data[883,0] = "'01483"
$data[884,0] = "'01484"
$data[885,0] = "'01485"
$data[886,0] = "'01486"
$data[887,0] = "'01487"
$data[888,0] = "'01488"
$data[889,0] = "'01489"
$data[890,0] = "'01490"
$data[891,0] = "'01491"
$data[892,0] = "'01492"
$data[893,0] = "'01493"
$data[894,0] = "'01494"
$data[895,0] = "'01495"
$data[896,0] = "'01496"
$data[897,0] = "'01497"
$data[898,0] = "'01498"
$data[899,0] = "'01500"
$data[900,0] = "'01501"
$data[901,0] = "'01502"
$data[902,0] = "'01503"
$data[903,0] = "'01504"
$data[904,0] = "'01505"
$data[905,0] = "'01506"
$data[906,0] = "'01507"
$data[907,0] = "'01508"
$data[908,0] = "'01509"
$data[909,0] = "'01511"
$data[910,0] = "'01512"
$data[911,0] = "'01513"
$data[912,0] = "'01514"
$data[913,0] = "'01515"
$data[914,0] = "'01516"
$data[915,0] = "'01517"
$data[916,0] = "'01518"
$data[917,0] = "'01519"
$data[918,0] = "'01520"
$data[919,0] = "'01521"
$data[920,0] = "'01522"
$data[921,0] = "'01523"
$data[922,0] = "'01524"
$data[923,0] = "'01525"
$data[924,0] = "'01526"
$data[925,0] = "'01527"
$data[926,0] = "'01528"
$data[927,0] = "'01529"
$data[928,0] = "'01530"
$data[929,0] = "'01531"
$data[930,0] = "'01532"
$data[931,0] = "'01533"
$data[932,0] = "'01534"
$data[933,0] = "'01535"
$data[934,0] = "'01536"
$data[935,0] = "'01537"
$data[936,0] = "'01538"
$data[937,0] = "'01539"
$data[938,0] = "'01540"
$data[939,0] = "'01541"
$data[940,0] = "'01542"
$data[941,0] = "'01543"
$data[942,0] = "'01544"
$data[943,0] = "'01545"
$data[944,0] = "'01546"
$data[945,0] = "'01547"
$data[946,0] = "'01548"
$data[947,0] = "'01549"
$data[948,0] = "'01550"
$data[949,0] = "'01551"
$data[950,0] = "'01552"
$data[951,0] = "'01553"
$data[952,0] = "'01554"
$data[953,0] = "'01555"
$data[954,0] = "'01556"
$data[955,0] = "'01557"
$data[956,0] = 123123
$data[957,0] = "'Dolore eius voluptat"
$data[958,0] = "'X1014"
$data[959,0] = "'X1023"
$data[960,0] = "'X1086"
$data[961,0] = "'X1094"
$data[962,0] = "'X1211"
$data[963,0] = "'X1212"
$data[964,0] = "'X1220"
$data[965,0] = "'X1224"
$data[966,0] = "'X1231"
$data[967,0] = "'X1235"
$data[968,0] = "'X1244"
$data[969,0] = "'X1250"
$data[970,0] = "'X1252"
$data[971,0] = "'X1257"
$data[972,0] = "'X1258"
$data[973,0] = "'X1259"
$data[974,0] = "'X1260"
$data[975,0] = "'X1261"
$data[976,0] = "'X1263"
$data[977,0] = "'X1264"
$data[978,0] = "'X1265"
$data[979,0] = "'X1266"
$data[980,0] = "'X1267"
$data[981,0] = "'X1268"
$data[982,0] = "'X1269"
$data[983,0] = "'X1270"
$data[984,0] = "'X1271"
$data[985,0] = "'X1272"
$data[986,0] = "'X1273"
$data[987,0] = "'X1274"
$data[988,0] = "'X1275"
$data[989,0] = "'X1276"
$data[990,0] = "'X1277"
$data[991,0] = "'X1278"
$data[992,0] = "'X1279"
$data[993,0] = "'X1280"
$data[994,0] = "'X1281"
$data[995,0] = "'X1282"
$data[996,0] = "'X1283"
$data[997,0] = "'X1284"
$data[998,0] = "'X1285"
$data[999,0] = "'X1286"
$data[1000,0] = "'X1287"
$data[1001,0] = "'X1288"
$data[1002,0] = "'X1289"
$data[1003,0] = "'X1291"
$data[1004,0] = "'X1292"
$data[1005,0] = "'X1293"
$data[1006,0] = "'X1294"
$data[1007,0] = "'X1295"
$data[1008,0] = "'X1296"
$data[1009,0] = "'X1297"
$data[1010,0] = "'X1298"
$data[1011,0] = "'X1299"
$data[1012,0] = "'X1300"
$data[1013,0] = "'X1301"
$data[1014,0] = "'X1304"
$data[1015,0] = "'X1306"
$data[1016,0] = "'X1307"
$data[1017,0] = "'X1308"
$data[1018,0] = "'X1309"
$data[1019,0] = "'X1310"
$data[1020,0] = "'X1311"
$data[1021,0] = "'X1312"
$data[1022,0] = "'X1313"
$data[1023,0] = "'X1314"
$data[1024,0] = "'X1316"
$data[1025,0] = "'X1317"
$data[1026,0] = "'X1318"
$data[1027,0] = "'X1319"
$data[1028,0] = "'X1320"
$data[1029,0] = "'X1321"
$data[1030,0] = "'X1322"
$data[1031,0] = "'X1323"
$data[1032,0] = "'X1324"
$data[1033,0] = "'X1325"
$data[1034,0] = "'X1326"
$data[1035,0] = "'X1327"
$data[1036,0] = "'X1328"
$data[1037,0] = "'X1329"
$data[1038,0] = "'X1330"
$data[1039,0] = "'X1331"
$data[1040,0] = "'X1332"
$data[1041,0] = "'X1333"
$data[1042,0] = "'X1334"
$data[1043,0] = "'X1335"
$data[1044,0] = "'X1336"
$data[1045,0] = "'X1337"
$data[1046,0] = "'X1338"
$data[1047,0] = "'X1339"
$data[1048,0] = "'X1340"
$data[1049,0] = "'X1341"
$data[1050,0] = "'X1342"
$data[1051,0] = "'X1343"
$data[1052,0] = "'X1344"
$data[1053,0] = "'X1345"
$data[1054,0] = "'X1346"
$data[1055,0] = "'X1347"
$data[1056,0] = "'X1348"
$data[1057,0] = "'X1349"
$data[1058,0] = "'X1350"
$data[1059,0] = "'X1351"
$data[1060,0] = "'X1352"
$data[1061,0] = "'X1353"
$data[1062,0] = "'X1354"
$data[1063,0] = "'X1355"
$data[1064,0] = "'X1356"
$data[1065,0] = "'X1357"
$data[1066,0] = "'X1358"
$data[1067,0] = "'X1359"
$data[1068,0] = "'X1360"
$data[1069,0] = "'X1361"
$data[1070,0] = "'X1362"
$data[1071,0] = "'X1363"
$data[1072,0] = "'X1364"
$data[1073,0] = "'X1365"
$data[1074,0] = "'X1366"
$data[1075,0] = "'X1367"
$data[1076,0] = "'X1368"
$data[1077,0] = "'X1369"
$data[1078,0] = "'X1370"
$data[1079,0] = "'X1371"
$data[1080,0] = "'X1372"
$data[1081,0] = "'X1373"
$data[1082,0] = "'X1374"
$data[1083,0] = "'X1375"
$data[1084,0] = "'X1376"
$data[1085,0] = "'X1378"
$data[1086,0] = "'X1379"
$data[1087,0] = "'X1380"
$data[1088,0] = "'X1381"
$data[1089,0] = "'X1382"
$data[1090,0] = "'X1383"
$data[1091,0] = "'X1384"
$data[1092,0] = "'X1385"
$data[1093,0] = "'X1386"
$data[1094,0] = "'X1387"
$data[1095,0] = "'X1388"
$data[1096,0] = "'X1389"
$data[1097,0] = "'X1390"
$data[1098,0] = "'X1391"
$data[1099,0] = "'X1392"
$data[1100,0] = "'X1393"
$data[1101,0] = "'X1394"
$data[1102,0] = "'X1395"
$data[1103,0] = "'X1396"
$data[1104,0] = "'X1397"
$data[1105,0] = "'X1398"
$data[1106,0] = "'X1399"
$data[1107,0] = "'X1400"
$data[1108,0] = "'X1401"
$data[1109,0] = "'X1402"
$data[1110,0] = "'X1403"
$data[1111,0] = "'X1404"
$data[1112,0] = "'X1405"
$data[1113,0] = "'X1406"
$data[1114,0] = "'X1407"
$data[1115,0] = "'X1408"
$data[1116,0] = "'X1409"
$data[1117,0] = "'X1410"
$data[1118,0] = "'X1411"
$data[1119,0] = "'X1412"
$data[1120,0] = "'X1413"
$data[1121,0] = "'X1414"
$data[1122,0] = "'X1415"
$data[1123,0] = "'X1416"
$data[1124,0] = "'X1417"
$data[1125,0] = "'X1418"
$data[1126,0] = "'X1419"
$data[1127,0] = "'X1420"
$data[1128,0] = "'X1421"
$data[1129,0] = "'X1422"
$data[1130,0] = "'XX006"
$data[1131,0] = "'XX016"
$data[1132,0] = "'XX018"
$data[1133,0] = "'XX022"
$data[1134,0] = "'XX026"
$data[1135,0] = "'XX027"
$data[1136,0] = "'XX034"
$data[1137,0] = "'XX158"
$data[1138,0] = "'XX159"
$data[1139,0] = "'XX355"
$data[1140,0] = "'XX356"
$data[1141,0] = "'XX437"
$data[1142,0] = "'XX439"
$data[1143,0] = "'XX804"
$data[1144,0] = "'XX845"
$data[1145,0] = "'XX866"
$data[1146,0] = "'XX969"
$data[1147,0] = "'XX970"
$data[1148,0] = "'ZZ003"
$data[1149,0] = "'ZZ006"
$data[1150,0] = "'ZZ009"
$data[1151,0] = "'ZZ021"
$data[1152,0] = "'ZZ022"
$data[1153,0] = "'ZZ026"
$data[1154,0] = "'ZZ027"
$data[1155,0] = "'ZZ028"
$data[1156,0] = "'ZZ029"
$data[1157,0] = "'ZZ030"
$data[1158,0] = "'ZZ031"
$data[1159,0] = "'ZZ032"

$target = $ws.Range("A2:A" + $lastNewRow)
$target.Value = $data

$ws.Range("A1:A1048576").Select()
